$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "3-132-873"
$ws.Range("C11").Value = "OGN FUSEHOLDER 5X20, 22.5 MM PITCH"
$ws.Range("H11").Value = "SCHURTER Inc."
